$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 12.85
$ws.Range("E7").Value = 13.089
$ws.Range("C8").Value = -12.378
$ws.Range("A12").Value = -21.882
$ws.Range("C12").Value = -13.002
$ws.Range("C14").Value = -11.675
$ws.Range("E19").Value = 12.513
$ws.Range("E21").Value = 13.147
$ws.Range("C22").Value = -12.473
$ws.Range("E24").Value = 12.848
